# Update gh-pages to output generated at 456a3b4
# ---------------------------------------------------------------------------
# This script reproduces, via the Excel COM object model, the data refresh
# that the commit applied to 广州-漫展信息.xlsx:
#   * a handful of "想去人数" (F column) counters were bumped on sheets
#     展览 (1), 演出 (2), 本地生活 (3) and 全部类型 (4)
#   * a brand-new event ("广州·司南 2024「出发」个人巡演·生日专场") was
#     inserted into 演出 (sheet 2) and 全部类型 (sheet 4), pushing the
#     following rows down by one and bumping their running index (col A)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force plain-text storage so Excel doesn't "helpfully" reinterpret
    # date-looking strings (e.g. "2024-05-04") as real dates, and make sure
    # we don't leave a stray number-format override behind afterwards.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-NumCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

function Set-IndexCell($ws, $row, $num) {
    # Column A holds a plain running index (row-2), styled like a header
    # label cell: bold, thin border, centered/top aligned.
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $num
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1
    $c.Borders.Weight = 2
}

function Set-EventRow($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i) {
    Set-TextCell $ws $row 2 $b
    Set-TextCell $ws $row 3 $c
    Set-TextCell $ws $row 4 $d
    Set-TextCell $ws $row 5 $e
    Set-NumCell  $ws $row 6 $f
    Set-NumCell  $ws $row 7 $g
    Set-TextCell $ws $row 8 $h
    Set-TextCell $ws $row 9 $i
}

# ===========================================================================
# Sheet 1: 展览 (Exhibition) — counter-only updates
# ===========================================================================
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    2  = 7761
    5  = 8400
    8  = 635
    9  = 448
    13 = 38
    14 = 75
    15 = 314
    17 = 264
    19 = 395
    20 = 148
    21 = 1086
    23 = 620
    24 = 2206
    25 = 737
    26 = 52
    29 = 614
    30 = 54
}
foreach ($row in $sheet1Updates.Keys) {
    Set-NumCell $ws1 $row 6 $sheet1Updates[$row]
}

# ===========================================================================
# Sheet 2: 演出 (Performance) — counter updates + new event row
# ===========================================================================
$ws2 = $wb.Worksheets.Item("演出")

Set-NumCell $ws2 2 6 288
Set-NumCell $ws2 4 6 327

# Insert a new row before the current last row (row 10); that row's existing
# content shifts down to row 11 automatically (formatting travels with it).
$ws2.Rows.Item(10).Insert()

# Row 10 is now blank — populate it with the new event, and restore the
# running-index cell (this row is still "item 9", same index as before).
Set-IndexCell $ws2 10 9
Set-EventRow $ws2 10 `
    "2024-05-04" `
    "广州·司南 2024「出发」个人巡演·生日专场" `
    "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)" `
    "2024.05.04 19:30-05.04 21:30" `
    0 `
    299 `
    "https://show.bilibili.com/platform/detail.html?id=83341" `
    "//i0.hdslb.com/bfs/openplatform/202403/hBiSFtLH1711365103423.jpeg"

# The old row 10 (now row 11) keeps its data/format, but its running index
# needs to advance from 9 to 10.
Set-NumCell $ws2 11 1 10

# ===========================================================================
# Sheet 3: 本地生活 (Local Life) — counter-only update
# ===========================================================================
$ws3 = $wb.Worksheets.Item("本地生活")
Set-NumCell $ws3 2 6 452

# ===========================================================================
# Sheet 4: 全部类型 (All Types) — counter updates + new event row
# ===========================================================================
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    2  = 452
    3  = 7762
    6  = 288
    7  = 8402
    10 = 635
    11 = 448
    15 = 327
    19 = 38
    20 = 75
    21 = 315
}
foreach ($row in $sheet4Updates.Keys) {
    Set-NumCell $ws4 $row 6 $sheet4Updates[$row]
}

# Insert the same new event before the old row 26 ("广州·第五人格only 1.0"),
# which (with everything after it) shifts down by one row.
$ws4.Rows.Item(26).Insert()

Set-IndexCell $ws4 26 25
Set-EventRow $ws4 26 `
    "2024-05-04" `
    "广州·司南 2024「出发」个人巡演·生日专场" `
    "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)" `
    "2024.05.04 19:30-05.04 21:30" `
    0 `
    299 `
    "https://show.bilibili.com/platform/detail.html?id=83341" `
    "//i0.hdslb.com/bfs/openplatform/202403/hBiSFtLH1711365103423.jpeg"

# Rows 27..41 are the old rows 26..40, shifted down one; bump their running
# index (col A) by one to keep it equal to row-2.
for ($row = 27; $row -le 41; $row++) {
    $newIndex = $row - 2
    Set-NumCell $ws4 $row 1 $newIndex
}

Write-Host "done"
